$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray leftover "Projects Summary" bullet list block
#    that follows "February 20, 2014" (paragraphs 4 through 15 in the
#    original document: the blank line, "Projects Summary" and its
#    bullet items, down to -- but not including -- the blank
#    paragraph that precedes "Summary").
# ------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(4)
$pEnd = $d.Paragraphs.Item(15)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from right before the "Kernel"
#    heading up to right after "February 20, 2014".
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Find the "February 20, 2014" run and collapse to its end.
$findRange = $d.Content
$ok = $findRange.Find.Execute("February 20, 2014", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Collapse(0)

# Work around collapsed-range placement quirks: insert a temporary
# marker, anchor the bookmark to it (non-collapsed), then remove the
# marker text -- the bookmark naturally stays collapsed in place.
$findRange.InsertAfter("@@TMPMARK@@")

$markRange = $d.Content
$ok = $markRange.Find.Execute("@@TMPMARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $markRange)

$removeRange = $d.Content
$ok = $removeRange.Find.Execute("@@TMPMARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ------------------------------------------------------------------
# 3) Drop the (now meaningless) <w:lastRenderedPageBreak/> cached in
#    the "Kernel" run by rewriting the paragraph's text in place.
# ------------------------------------------------------------------
$kernelRange = $d.Content
$ok = $kernelRange.Find.Execute("Kernel", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$kernelPara = $kernelRange.Paragraphs.Item(1)
$kernelPara.Range.Text = "Kernel"

Write-Host "edit complete"
